# Saldo.xlsx update — "Add files via upload"
#
# Net effect on the `Export` sheet's data table (Conta | Nome | Saldo):
#   + insert 005639338 | CLOTILDE | 20000   (right before the 004364200 BLOCO row)
#   - remove 004329030 | DANIELA  | 14040.07
#   - remove 004927044 | CINTIA   | 10000
#   - remove 001090818 | MARIA    | 9987.61
#   - remove 004460491 | PEDRO    | 3383.4   (replaced in-place by MARCELO/8000)
#   + replace that row with 004748761 | MARCELO | 8000
#   - remove the older 004748761 | MARCELO | 1000 row further down
#   - remove 004231371 | ADRIANO  | 350
#
# Work from the bottom of the sheet upward so row numbers used below stay
# valid as rows are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove 004231371 | ADRIANO | 350  (row 34)
$ws.Rows.Item(34).Delete()

# 2) Remove the original 004748761 | MARCELO | 1000 row (row 21)
$ws.Rows.Item(21).Delete()

# 3) Remove 004927044 | CINTIA | 10000  (row 12)
$ws.Rows.Item(12).Delete()

# 4) Remove 001090818 | MARIA | 9987.61  (row 12, after the previous delete)
$ws.Rows.Item(12).Delete()

# 5) The row that was 004460491 | PEDRO | 3383.4 is now row 12 — turn it into
#    004748761 | MARCELO | 8000 in place.
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value2 = "004748761"
$ws.Cells.Item(12, 1).ClearFormats()
$ws.Cells.Item(12, 2).Value = "MARCELO"
$ws.Cells.Item(12, 3).Value = 8000

# 6) Remove 004329030 | DANIELA | 14040.07  (row 10)
$ws.Rows.Item(10).Delete()

# 7) Insert the new 005639338 | CLOTILDE | 20000 row right before 004364200
#    (row 8).
$ws.Rows.Item(8).Insert()
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value2 = "005639338"
$ws.Cells.Item(8, 1).ClearFormats()
$ws.Cells.Item(8, 2).Value = "CLOTILDE"
$ws.Cells.Item(8, 3).Value = 20000
